# Set exercises Q1 - updated graphs axis, Q1C valgrind analysis

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "D1 misses" / "D refs" / "33% miss rate" (valgrind) columns ---
$ws.Range("H1").Value = "D1 misses"
$ws.Range("J1").Value = "D refs"
$ws.Range("L1").Value = "33% miss rate"

$ws.Range("H2").Value = 2123067
$ws.Range("J2").Value = 6381601

$ws.Range("H2:H3").NumberFormat = "#,##0"
$ws.Range("J2:J3").NumberFormat = "#,##0"

# --- Chart 1 ("a." gigaFLOPS vs N) -> "q1a.c GFLOPS vs N" ---
$chart1 = $ws.ChartObjects(1).Chart
$chart1.ChartTitle.Text = "q1a.c GFLOPS vs N (ubuntu virtual machine)"
$chart1.Axes(2).AxisTitle.Text = "GFLOPS"

# --- Chart 2 ("b." gigaFLOPS vs N) -> "q1b.c GFLOPS vs N" ---
$chart2 = $ws.ChartObjects(2).Chart
$chart2.ChartTitle.Text = "q1b.c GFLOPS vs N (ubuntu virtual machine)"
$chart2.Axes(2).AxisTitle.Text = "GFLOPS"

# --- Selection / view state ---
$ws.Range("K8").Select()
